$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.715.03'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '2.111.80'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.23'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.391'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0780'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '2.429.96'
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.787'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '2.124.42'
$ws.Range("E17").Value = '  +2.38%  '
$ws.Range("D18").Value = '37.748.61'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.135'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.119'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0621'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.57'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0965'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.00%  '
$ws.Range("D43").Value = '1.475.55'
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0214'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.96%  '
$ws.Range("D51").Value = '2.314.33'
$ws.Range("E51").Value = '  +2.01%  '
